# Append a new row (row 48) to each of the four data sheets. The new row
# duplicates the last existing row (row 47) except for column A, which
# receives a new timestamp value.

$wb = $excel.ActiveWorkbook

$newTimestamp = [double]"45834.49655092593"

$sheetsData = @(
    @{
        Index = 1
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x5C"
        E = "0xf"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 348
        I = 15
    },
    @{
        Index = 2
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x6C"
        E = "0xe"
        F = 400
        G = [double]"5.68432987514711e+23"
        H = 364
        I = 14
    },
    @{
        Index = 3
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x69"
        E = "0x3"
        F = 110
        G = [double]"5.68631262647114e+23"
        H = 105
        I = 3
    },
    @{
        Index = 4
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x69"
        E = "0x3"
        F = 110
        G = [double]"9.85046333984776e+23"
        H = 105
        I = 3
    }
)

foreach ($sheetData in $sheetsData) {
    $ws = $wb.Worksheets.Item($sheetData.Index)

    $row = 48

    $ws.Cells.Item($row, 1).Value = $newTimestamp
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 2).Value = $sheetData.B
    $ws.Cells.Item($row, 3).Value = $sheetData.C
    $ws.Cells.Item($row, 4).Value = $sheetData.D
    $ws.Cells.Item($row, 5).Value = $sheetData.E
    $ws.Cells.Item($row, 6).Value = $sheetData.F
    $ws.Cells.Item($row, 7).Value = $sheetData.G
    $ws.Cells.Item($row, 8).Value = $sheetData.H
    $ws.Cells.Item($row, 9).Value = $sheetData.I
}
